$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Cabernet Franc"
$ws.Range("B1").Value = 1300
$ws.Range("C1").Value = "Bodega Cinco"
$ws.Range("D1").Value = "Famoso por sus vinos tintos de alta calidad"
$ws.Range("E1").Value = "Argentina"
$ws.Range("F1").Value = "Cabernet FrancSauvignon Blanc es una variedad de uva blanca conocida por su frescura y sus aromas herbáceos y cítricos.Bonarda es una variedad de uva tinta que se utiliza en la producción de vinos tintos suaves y afrutados, con sabores a frutas negras y especias.Sangiovese es una variedad de uva tinta que se asocia principalmente con los vinos italianos, conocidos por su acidez y sabor a frutas rojas."
$ws.Range("G1").Value = 9.9

$ws.Range("A2").Value = "Malbec"
$ws.Range("B2").Value = 1200
$ws.Range("C2").Value = "Bodega Cinco"
$ws.Range("D2").Value = "Famoso por sus vinos tintos de alta calidad"
$ws.Range("E2").Value = "Argentina"
$ws.Range("F2").Value = "MalbecMalbec"
$ws.Range("G2").Value = 9.2

$ws.Range("A3").Value = "Torrontés Clásico"
$ws.Range("B3").Value = 1150
$ws.Range("C3").Value = "Bodega Dos"
$ws.Range("D3").Value = "Reconocida por sus Malbecs"
$ws.Range("E3").Value = "Argentina"
$ws.Range("F3").Value = "Torrontés ClásicoMourvèdre es una variedad de uva tinta que se utiliza en la producción de vinos tintos robustos y especiados.Cabernet Sauvignon es una variedad de uva tinta ampliamente reconocida por su presencia en los vinos tintos de Bordeaux.Garnacha Blanca es una variedad de uva blanca que produce vinos blancos con cuerpo y textura, con sabores a frutas blancas y notas florales."
$ws.Range("G3").Value = 9

$ws.Range("A4").Value = "Merlot Reserva"
$ws.Range("B4").Value = 1250
$ws.Range("C4").Value = "Bodega Tres"
$ws.Range("D4").Value = "Región importante de San Juan"
$ws.Range("E4").Value = "Argentina"
$ws.Range("F4").Value = "Merlot ReservaMerlot Reserva"
$ws.Range("G4").Value = 8.9

$ws.Range("A5").Value = "Riesling"
$ws.Range("B5").Value = 1350
$ws.Range("C5").Value = "Bodega Ocho"
$ws.Range("D5").Value = "Famoso por sus vinos tintos de alta calidad"
$ws.Range("E5").Value = "Argentina"
$ws.Range("F5").Value = "RieslingRiesling"
$ws.Range("G5").Value = 8.6

$ws.Range("A6").Value = "Chardonnay"
$ws.Range("B6").Value = 1350
$ws.Range("C6").Value = "Bodega Cuatro"
$ws.Range("D6").Value = "Famoso por sus vinos tintos de alta calidad"
$ws.Range("E6").Value = "Argentina"
$ws.Range("F6").Value = "Sémillon es una variedad de uva blanca que se utiliza en la producción de vinos blancos secos, dulces y también vinos de postre.Pinot Grigio es una variedad de uva blanca que produce vinos blancos ligeros y refrescantes, con notas cítricas y florales.Garnacha Blanca es una variedad de uva blanca que produce vinos blancos con cuerpo y textura, con sabores a frutas blancas y notas florales."
$ws.Range("G6").Value = 8.4

$ws.Range("A7").Value = "Zinfandel"
$ws.Range("B7").Value = 1100
$ws.Range("C7").Value = "Bodega Ocho"
$ws.Range("D7").Value = "Famoso por sus vinos tintos de alta calidad"
$ws.Range("E7").Value = "Argentina"
$ws.Range("F7").Value = "ZinfandelZinfandel"
$ws.Range("G7").Value = 8.325

$ws.Range("A8").Value = "Cabernet"
$ws.Range("B8").Value = 1300
$ws.Range("C8").Value = "Bodega Seis"
$ws.Range("D8").Value = "Reconocida por sus Malbecs"
$ws.Range("E8").Value = "Argentina"
$ws.Range("F8").Value = "CabernetMerlot es una variedad de uva tinta que se caracteriza por su suavidad y sabor frutal en los vinos.Pinot Noir es una variedad de uva tinta que se asocia con vinos ligeros, elegantes y afrutados."
$ws.Range("G8").Value = 8.3

$ws.Range("A9").Value = "Pinot Noir"
$ws.Range("B9").Value = 1500
$ws.Range("C9").Value = "Bodega Tres"
$ws.Range("D9").Value = "Región importante de San Juan"
$ws.Range("E9").Value = "Argentina"
$ws.Range("F9").Value = "Garnacha Blanca es una variedad de uva blanca que produce vinos blancos con cuerpo y textura, con sabores a frutas blancas y notas florales.Tannat es una variedad de uva tinta que se asocia principalmente con los vinos de Uruguay, conocidos por su estructura tánica y sabor intenso."
$ws.Range("G9").Value = 8.25

$ws.Range("A10").Value = "Torrontés"
$ws.Range("B10").Value = 1150
$ws.Range("C10").Value = "Bodega Ocho"
$ws.Range("D10").Value = "Famoso por sus vinos tintos de alta calidad"
$ws.Range("E10").Value = "Argentina"
$ws.Range("F10").Value = "TorrontésRiesling es una variedad de uva blanca que puede producir desde vinos secos y refrescantes hasta vinos dulces y aromáticos."
$ws.Range("G10").Value = 8.2
